# Update "想去人数" (want-to-go count) values in column F on the
# "展览" and "全部类型" worksheets to reflect refreshed data as of 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 922
$ws1.Range("F3").Value = 550
$ws1.Range("F7").Value = 321
$ws1.Range("F11").Value = 186
$ws1.Range("F12").Value = 4665
$ws1.Range("F15").Value = 465
$ws1.Range("F18").Value = 314
$ws1.Range("F22").Value = 693
$ws1.Range("F24").Value = 288
$ws1.Range("F27").Value = 1672

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 922
$ws4.Range("F7").Value = 550
$ws4.Range("F12").Value = 321
$ws4.Range("F16").Value = 186
$ws4.Range("F17").Value = 186
$ws4.Range("F18").Value = 4665
$ws4.Range("F22").Value = 465
$ws4.Range("F25").Value = 314
$ws4.Range("F33").Value = 693
$ws4.Range("F38").Value = 288
$ws4.Range("F41").Value = 1672
